$d = $word.ActiveDocument

# --- Change 1: update the "Curso (semestre ideal)" line ---
$d.Content.Find.Execute(
    "Curso (semestre ideal): EA (4), EB (3), EQD (3), EQN (4)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EA (4), EB (5), EQD (4), EQN (5)", 2)

# --- Change 2: split the "LOB1004" requisito bullet into two bullets that
#     live inside the same paragraph (separated by a manual line break) ---

# 2a. Rename the existing requisito text in place (keeps the run + its <w:br/>).
$d.Content.Find.Execute(
    "LOB1004 -  Cálculo II  (Requisito fraco)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "LOB1024 -  Mecânica  (Requisito fraco)", 2)

# 2b. Append a brand-new paragraph right after it (same ListBullet style) and
#     give it the second requisito's text, ending with a manual line break
#     character so Word records a trailing <w:br/>, matching the first run.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1024*") {
        $p.Range.InsertParagraphAfter()
    }
}
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count)
$manualBreak = [char]11
$newPara.Range.Text = "LOB1052 -  Cálculo III  (Requisito fraco)" + $manualBreak

# 2c. Merge the new paragraph back into the previous one by deleting the
#     paragraph mark between them - this leaves a single <w:p> containing
#     two <w:r> runs (one per requisito), each ending in its own <w:br/>.
$firstPara = $d.Paragraphs.Item($count - 1)
$markPos = $firstPara.Range.End - 1
$d.Range($markPos, $markPos + 1).Delete()
